# Scheduled-runner refresh of market-price-derived Leve profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all eight job
# sheets, per the latest price pull. Only H..N columns change; A..G
# (leve metadata / item ids) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1786.8334
$ws.Range("J17").Value = 1786.8334
$ws.Range("L17").Value = 5360.5002
$ws.Range("N17").Value = -5696.5002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3678.8333
$ws.Range("I76").Value = 3797
$ws.Range("K76").Value = 3797
$ws.Range("M76").Value = -3482

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3678.8333
$ws.Range("I79").Value = 3797
$ws.Range("K79").Value = 3797
$ws.Range("M79").Value = -2705

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 75000
$ws.Range("I86").Value = 75000
$ws.Range("K86").Value = 75000
$ws.Range("M86").Value = -73877

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 75000
$ws.Range("I89").Value = 75000
$ws.Range("K89").Value = 375000
$ws.Range("M89").Value = -369384

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2646
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2646
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 7938
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -10154

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 747.25
$ws.Range("I118").Value = 747.25
$ws.Range("K118").Value = 2241.75
$ws.Range("M118").Value = -584.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1675.909
$ws.Range("J127").Value = 1113.5714
$ws.Range("L127").Value = 3340.7142
$ws.Range("N127").Value = -13260.7142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1994.6
$ws.Range("I132").Value = 2025.3846
$ws.Range("J132").Value = 1794.5
$ws.Range("K132").Value = 6076.1538
$ws.Range("L132").Value = 5383.5
$ws.Range("M132").Value = -3546.1538
$ws.Range("N132").Value = -10443.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1893.2
$ws.Range("I102").Value = 1893.2
$ws.Range("K102").Value = 1893.2
$ws.Range("M102").Value = -271.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 425130.5
$ws.Range("I133").Value = 350000
$ws.Range("K133").Value = 350000
$ws.Range("M133").Value = -347470

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 35000
$ws.Range("J45").Value = 35000
$ws.Range("L45").Value = 35000
$ws.Range("N45").Value = -36616

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 4999
$ws.Range("I46").Value = 4999
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4999
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -4701
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 46547.332
$ws.Range("J117").Value = 46547.332
$ws.Range("L117").Value = 46547.332
$ws.Range("N117").Value = -55725.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 99999.2
$ws.Range("J138").Value = 99999.2
$ws.Range("L138").Value = 99999.2
$ws.Range("N138").Value = -110279.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 99.46154
$ws.Range("I7").Value = 74.875
$ws.Range("K7").Value = 74.875
$ws.Range("M7").Value = 38.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 623.75
$ws.Range("I22").Value = 623.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 623.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -273.75
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1145.25
$ws.Range("I31").Value = 1145.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1145.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -850.25
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1145.25
$ws.Range("I34").Value = 1145.25
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1145.25
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -943.25
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2647.4348
$ws.Range("I132").Value = 2811.6667
$ws.Range("K132").Value = 8435.000100000001
$ws.Range("M132").Value = -5905.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 293.86667
$ws.Range("I12").Value = 301.81818
$ws.Range("J12").Value = 272
$ws.Range("K12").Value = 905.45454
$ws.Range("L12").Value = 816
$ws.Range("M12").Value = -732.45454
$ws.Range("N12").Value = -1162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 6666
$ws.Range("I125").Value = 6666
$ws.Range("K125").Value = 19998
$ws.Range("M125").Value = -15078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 12341.143
$ws.Range("I140").Value = 1045.25
$ws.Range("K140").Value = 3135.75
$ws.Range("M140").Value = 2044.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1893328.6
$ws.Range("I11").Value = 2042216.6
$ws.Range("K11").Value = 2042216.6
$ws.Range("M11").Value = -2042077.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 50325
$ws.Range("J136").Value = 50325
$ws.Range("L136").Value = 150975
$ws.Range("N136").Value = -156075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 770.6667
$ws.Range("I16").Value = 770.6667
$ws.Range("K16").Value = 770.6667
$ws.Range("M16").Value = -600.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 12500
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 1750
$ws.Range("J46").Value = 4200
$ws.Range("K46").Value = 1750
$ws.Range("L46").Value = 4200
$ws.Range("M46").Value = -1562
$ws.Range("N46").Value = -4576

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2294.3333
$ws.Range("I61").Value = 2206.125
$ws.Range("K61").Value = 2206.125
$ws.Range("M61").Value = -2004.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2294.3333
$ws.Range("I113").Value = 2206.125
$ws.Range("K113").Value = 2206.125
$ws.Range("M113").Value = -36.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5127
$ws.Range("I132").Value = 5720.4287
$ws.Range("J132").Value = 3742.3333
$ws.Range("K132").Value = 17161.2861
$ws.Range("L132").Value = 11226.9999
$ws.Range("M132").Value = -14631.2861
$ws.Range("N132").Value = -16286.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 32000
$ws.Range("J116").Value = 32000
$ws.Range("L116").Value = 32000
$ws.Range("N116").Value = -41178

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1661.2778
$ws.Range("I132").Value = 1641.0588
$ws.Range("J132").Value = 2005
$ws.Range("K132").Value = 4923.1764
$ws.Range("L132").Value = 6015
$ws.Range("M132").Value = -2393.1764
$ws.Range("N132").Value = -11075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 525107.5
$ws.Range("J135").Value = 525107.5
$ws.Range("L135").Value = 525107.5
$ws.Range("N135").Value = -535247.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 54000
$ws.Range("J137").Value = 54000
$ws.Range("L137").Value = 54000
$ws.Range("N137").Value = -64200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279
